# Update values on the "rawdata" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rawdata")

$ws.Range("B2").Value = 4600
$ws.Range("C2").Value = 56
$ws.Range("D2").Value = 36

$ws.Range("B3").Value = 4260
$ws.Range("B4").Value = 4280
$ws.Range("B5").Value = 4290
$ws.Range("B6").Value = 4320
$ws.Range("B7").Value = 4400

$ws.Range("B8").Value = 4500
$ws.Range("C8").Value = 57
$ws.Range("D8").Value = 38

# Move the active selection from D8 to D9 as reflected in the diff
$ws.Activate()
$ws.Range("D9").Select()
